$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cmin row (row 5): B5 1000 -> 900, C5 0 -> 500
$ws.Range("B5").Value = 900
$ws.Range("C5").Value = 500

# cmax row (row 6): C6 3500 -> 6500
$ws.Range("C6").Value = 6500

# ofthresmin row (row 15): B15 0.1 -> 2, C15 0.1 -> 2
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 2

# ofthresmax row (row 16): B16 10 -> 10000, C16 10 -> 10000
$ws.Range("B16").Value = 10000
$ws.Range("C16").Value = 10000

# Update selection to C6 as active cell
$ws.Range("C6").Select()
